$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the hours and summary for the week of row 14 (week starting 45488)
$ws.Range("B14").Value = "6, 45"
$ws.Range("C14").Value = "preprocessing session 4"

# Move the active selection to C15, matching the end-of-session cursor position
$ws.Range("C15").Select()
